$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new row before the current row 19 (the "Commands To Imaging Program"
# section header), which pushes the whole second table down by one row
# (old rows 19-34 become rows 20-35).
$ws.Rows("19:19").Insert()

# Fill the two newly freed rows (17 and 18) at the end of the first table
# ("Commands To SpineTracker") with the new commands.
$ws.Range("A17").Value2 = "CustomCommandReceived"

$ws.Range("A18").Value2 = "UncagingLocation"
$ws.Range("B18").Value2 = "ROI_x_pixels, ROI_y_pixels"
$ws.Range("D18").Value2 = "sent as an answer to SetUncagingLocation"

# Append the two new rows to the end of the second table
# ("Commands To Imaging Program").
$ws.Range("A36").Value2 = "CustomCommand"
$ws.Range("B36").Value2 = "string"
$ws.Range("C36").Value2 = "CustomCommand,page_acq"
$ws.Range("D36").Value2 = "CustomCommandReceived"
$ws.Range("E36").Value2 = "send a custom command to imaging program immediately before executing a step"

$ws.Range("A37").Value2 = "SetUncagingLocation"
$ws.Range("B37").Value2 = "ROI_x_pixels, ROI_y_pixels"
$ws.Range("C37").Value2 = "SetUncagingLocation,37,42"
$ws.Range("D37").Value2 = "UncagingLocation"
$ws.Range("E37").Value2 = "sends an uncaging location without actually uncaging"

# Refresh the remembered sort-state on the now-shifted "Commands To
# Imaging Program" table so it still points at the right block.
$so = $ws.Sort
$so.SortFields.Clear()
$so.SortFields.Add($ws.Range("A20"))
$so.SetRange($ws.Range("A21:E35"))
$so.Header = 2
$so.Apply()

# Update the selection to match the edited cell.
$ws.Range("D18").Select()
